$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028387484464507
$ws.Cells.Item(2, 4).Value = 1.031031937333478
$ws.Cells.Item(2, 5).Value = 1.037024270101197
$ws.Cells.Item(2, 6).Value = 1.044831663553818
$ws.Cells.Item(2, 9).Value = 1.031458094441068
$ws.Cells.Item(2, 10).Value = 1.033540083901993
$ws.Cells.Item(2, 11).Value = 1.033841229938427
$ws.Cells.Item(2, 12).Value = 1.039816336229628
$ws.Cells.Item(2, 13).Value = 1.047601607058841
$ws.Cells.Item(2, 14).Value = 1.015198144459557
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029329010817944
$ws.Cells.Item(3, 4).Value = 1.031886700390661
$ws.Cells.Item(3, 5).Value = 1.037895119148667
$ws.Cells.Item(3, 6).Value = 1.045867775930881
$ws.Cells.Item(3, 9).Value = 1.031583434463904
$ws.Cells.Item(3, 10).Value = 1.034122238610532
$ws.Cells.Item(3, 11).Value = 1.034504347739644
$ws.Cells.Item(3, 12).Value = 1.040496745678366
$ws.Cells.Item(3, 13).Value = 1.04844845112498
$ws.Cells.Item(3, 14).Value = 1.015390744270464
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.029938115781053
$ws.Cells.Item(4, 4).Value = 1.03243993731929
$ws.Cells.Item(4, 5).Value = 1.03845896593425
$ws.Cells.Item(4, 6).Value = 1.046538867688036
$ws.Cells.Item(4, 9).Value = 1.031662283087776
$ws.Cells.Item(4, 10).Value = 1.03449820466779
$ws.Cells.Item(4, 11).Value = 1.034932938515112
$ws.Cells.Item(4, 12).Value = 1.040936728023079
$ws.Cells.Item(4, 13).Value = 1.048996469606627
$ws.Cells.Item(4, 14).Value = 1.015515108448591
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030194152657535
$ws.Cells.Item(5, 4).Value = 1.032672552384074
$ws.Cells.Item(5, 5).Value = 1.038696089599414
$ws.Cells.Item(5, 6).Value = 1.046821151123285
$ws.Cells.Item(5, 9).Value = 1.031694890465212
$ws.Cells.Item(5, 10).Value = 1.034656085857271
$ws.Cells.Item(5, 11).Value = 1.03511299958756
$ws.Cells.Item(5, 12).Value = 1.041121626709795
$ws.Cells.Item(5, 13).Value = 1.049226868476703
$ws.Cells.Item(5, 14).Value = 1.01556732839164
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030237140543461
$ws.Cells.Item(6, 4).Value = 1.032711611451293
$ws.Cells.Item(6, 5).Value = 1.038735908509744
$ws.Cells.Item(6, 6).Value = 1.046868556899635
$ws.Cells.Item(6, 9).Value = 1.031700333664033
$ws.Cells.Item(6, 10).Value = 1.034682584518644
$ws.Cells.Item(6, 11).Value = 1.03514322565656
$ws.Cells.Item(6, 12).Value = 1.041152667887358
$ws.Cells.Item(6, 13).Value = 1.049265554130441
$ws.Cells.Item(6, 14).Value = 1.015576092660183
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.029941537080554
$ws.Cells.Item(7, 4).Value = 1.032443045400062
$ws.Cells.Item(7, 5).Value = 1.038462134069623
$ws.Cells.Item(7, 6).Value = 1.046542638959609
$ws.Cells.Item(7, 9).Value = 1.031662720914289
$ws.Cells.Item(7, 10).Value = 1.034500314972521
$ws.Cells.Item(7, 11).Value = 1.034935344966825
$ws.Cells.Item(7, 12).Value = 1.040939198924183
$ws.Cells.Item(7, 13).Value = 1.048999548162071
$ws.Cells.Item(7, 14).Value = 1.015515806461036
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028705703223361
$ws.Cells.Item(8, 4).Value = 1.031320777138377
$ws.Cells.Item(8, 5).Value = 1.037318504629932
$ws.Cells.Item(8, 6).Value = 1.045181685781894
$ws.Cells.Item(8, 9).Value = 1.031500920076984
$ws.Cells.Item(8, 10).Value = 1.033736975793726
$ws.Cells.Item(8, 11).Value = 1.034065434728197
$ws.Cells.Item(8, 12).Value = 1.040046342878873
$ws.Cells.Item(8, 13).Value = 1.047887790268563
$ws.Cells.Item(8, 14).Value = 1.015263288267741
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026527078928655
$ws.Cells.Item(9, 4).Value = 1.029344383040099
$ws.Cells.Item(9, 5).Value = 1.035306007274238
$ws.Cells.Item(9, 6).Value = 1.042788594123732
$ws.Cells.Item(9, 9).Value = 1.031198571993396
$ws.Cells.Item(9, 10).Value = 1.03238634624392
$ws.Cells.Item(9, 11).Value = 1.032528826259711
$ws.Cells.Item(9, 12).Value = 1.038470850608647
$ws.Cells.Item(9, 13).Value = 1.04592918450924
$ws.Cells.Item(9, 14).Value = 1.014816335923118
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025074083152964
$ws.Cells.Item(10, 4).Value = 1.028027652403366
$ws.Cells.Item(10, 5).Value = 1.033966238250516
$ws.Cells.Item(10, 6).Value = 1.041196673446645
$ws.Cells.Item(10, 9).Value = 1.030985464365345
$ws.Cells.Item(10, 10).Value = 1.031482260848811
$ws.Cells.Item(10, 11).Value = 1.031501979170597
$ws.Cells.Item(10, 12).Value = 1.037419118987907
$ws.Cells.Item(10, 13).Value = 1.044623807092624
$ws.Cells.Item(10, 14).Value = 1.014517053048315
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024444790880234
$ws.Cells.Item(11, 4).Value = 1.027457712523307
$ws.Cells.Item(11, 5).Value = 1.033386566181675
$ws.Cells.Item(11, 6).Value = 1.040508189478023
$ws.Cells.Item(11, 9).Value = 1.030890457885435
$ws.Cells.Item(11, 10).Value = 1.031089922708757
$ws.Cells.Item(11, 11).Value = 1.031056775349211
$ws.Cells.Item(11, 12).Value = 1.036963385287218
$ws.Cells.Item(11, 13).Value = 1.044058660490133
$ws.Cells.Item(11, 14).Value = 1.014387152197579
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024211023933834
$ws.Cells.Item(12, 4).Value = 1.027246044547693
$ws.Cells.Item(12, 5).Value = 1.033171319778631
$ws.Cells.Item(12, 6).Value = 1.040252581141893
$ws.Cells.Item(12, 9).Value = 1.030854758853289
$ws.Cells.Item(12, 10).Value = 1.030944061930984
$ws.Cells.Item(12, 11).Value = 1.030891321704181
$ws.Cells.Item(12, 12).Value = 1.036794057165835
$ws.Cells.Item(12, 13).Value = 1.043848754272705
$ws.Cells.Item(12, 14).Value = 1.014338855036598
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024261168623042
$ws.Cells.Item(13, 4).Value = 1.027291446533696
$ws.Cells.Item(13, 5).Value = 1.033217487696786
$ws.Cells.Item(13, 6).Value = 1.040307404315786
$ws.Cells.Item(13, 9).Value = 1.030862434932104
$ws.Cells.Item(13, 10).Value = 1.030975355384868
$ws.Cells.Item(13, 11).Value = 1.03092681591109
$ws.Cells.Item(13, 12).Value = 1.036830380797245
$ws.Cells.Item(13, 13).Value = 1.043893779204634
$ws.Cells.Item(13, 14).Value = 1.014349217027789
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024425468031617
$ws.Cells.Item(14, 4).Value = 1.027440215301728
$ws.Cells.Item(14, 5).Value = 1.033368772421061
$ws.Cells.Item(14, 6).Value = 1.040487058253031
$ws.Cells.Item(14, 9).Value = 1.030887515339962
$ws.Cells.Item(14, 10).Value = 1.031077868436988
$ws.Cells.Item(14, 11).Value = 1.031043100630885
$ws.Cells.Item(14, 12).Value = 1.036949389553323
$ws.Cells.Item(14, 13).Value = 1.04404130927844
$ws.Cells.Item(14, 14).Value = 1.014383160880058
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024526695716479
$ws.Cells.Item(15, 4).Value = 1.027531881063244
$ws.Cells.Item(15, 5).Value = 1.033461993186644
$ws.Cells.Item(15, 6).Value = 1.040597765596857
$ws.Cells.Item(15, 9).Value = 1.030902913978226
$ws.Cells.Item(15, 10).Value = 1.031141013050779
$ws.Cells.Item(15, 11).Value = 1.031114736230345
$ws.Cells.Item(15, 12).Value = 1.0370227083946
$ws.Cells.Item(15, 13).Value = 1.044132209357397
$ws.Cells.Item(15, 14).Value = 1.014404068693994
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.0251158444033
$ws.Cells.Item(16, 4).Value = 1.028065482016923
$ws.Cells.Item(16, 5).Value = 1.034004718884711
$ws.Cells.Item(16, 6).Value = 1.041242383380871
$ws.Cells.Item(16, 9).Value = 1.030991712183697
$ws.Cells.Item(16, 10).Value = 1.031508280917668
$ws.Cells.Item(16, 11).Value = 1.03153151389811
$ws.Cells.Item(16, 12).Value = 1.037449357711429
$ws.Cells.Item(16, 13).Value = 1.04466131598574
$ws.Cells.Item(16, 14).Value = 1.014525667645523
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.02548536594166
$ws.Cells.Item(17, 4).Value = 1.028400253611286
$ws.Cells.Item(17, 5).Value = 1.03434527936542
$ws.Cells.Item(17, 6).Value = 1.041646957409848
$ws.Cells.Item(17, 9).Value = 1.031046682676032
$ws.Cells.Item(17, 10).Value = 1.031738427785412
$ws.Cells.Item(17, 11).Value = 1.031792794912934
$ws.Cells.Item(17, 12).Value = 1.037716896577057
$ws.Cells.Item(17, 13).Value = 1.044993235378394
$ws.Cells.Item(17, 14).Value = 1.014601860821943
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.025700888489385
$ws.Cells.Item(18, 4).Value = 1.028595540717061
$ws.Cells.Item(18, 5).Value = 1.034543966536677
$ws.Cells.Item(18, 6).Value = 1.041883018448495
$ws.Cells.Item(18, 9).Value = 1.031078482635009
$ws.Cells.Item(18, 10).Value = 1.03187258522286
$ws.Cells.Item(18, 11).Value = 1.031945140437326
$ws.Cells.Item(18, 12).Value = 1.037872915870675
$ws.Cells.Item(18, 13).Value = 1.0451868470162
$ws.Cells.Item(18, 14).Value = 1.014646273142211
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025774373866423
$ws.Cells.Item(19, 4).Value = 1.028662132018455
$ws.Cells.Item(19, 5).Value = 1.034611721161258
$ws.Cells.Item(19, 6).Value = 1.041963522681065
$ws.Cells.Item(19, 9).Value = 1.03108928090143
$ws.Cells.Item(19, 10).Value = 1.031918315283432
$ws.Cells.Item(19, 11).Value = 1.031997076923217
$ws.Cells.Item(19, 12).Value = 1.037926109033804
$ws.Cells.Item(19, 13).Value = 1.045252865027533
$ws.Cells.Item(19, 14).Value = 1.014661411520558
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.025445721129816
$ws.Cells.Item(20, 4).Value = 1.02836433365955
$ws.Cells.Item(20, 5).Value = 1.03430873589003
$ws.Cells.Item(20, 6).Value = 1.041603542179731
$ws.Cells.Item(20, 9).Value = 1.031040812102107
$ws.Cells.Item(20, 10).Value = 1.03171374381892
$ws.Cells.Item(20, 11).Value = 1.031764767643492
$ws.Cells.Item(20, 12).Value = 1.037688195451525
$ws.Cells.Item(20, 13).Value = 1.044957622665943
$ws.Cells.Item(20, 14).Value = 1.014593689102783
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024377086492142
$ws.Cells.Item(21, 4).Value = 1.027396405686193
$ws.Cells.Item(21, 5).Value = 1.033324220910851
$ws.Cells.Item(21, 6).Value = 1.040434151181647
$ws.Cells.Item(21, 9).Value = 1.030880141080888
$ws.Cells.Item(21, 10).Value = 1.031047684451053
$ws.Cells.Item(21, 11).Value = 1.031008860022237
$ws.Cells.Item(21, 12).Value = 1.036914345768801
$ws.Cells.Item(21, 13).Value = 1.043997864944999
$ws.Cells.Item(21, 14).Value = 1.014373166534726
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023705079341809
$ws.Cells.Item(22, 4).Value = 1.026788022810654
$ws.Cells.Item(22, 5).Value = 1.032705620514006
$ws.Cells.Item(22, 6).Value = 1.03969963380159
$ws.Cells.Item(22, 9).Value = 1.030776752525817
$ws.Cells.Item(22, 10).Value = 1.030628160949615
$ws.Cells.Item(22, 11).Value = 1.030533098819067
$ws.Cells.Item(22, 12).Value = 1.036427516503434
$ws.Cells.Item(22, 13).Value = 1.043394510847176
$ws.Cells.Item(22, 14).Value = 1.014234248075235
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024061333673153
$ws.Cells.Item(23, 4).Value = 1.02711051955119
$ws.Cells.Item(23, 5).Value = 1.033033513736829
$ws.Cells.Item(23, 6).Value = 1.04008894634758
$ws.Cells.Item(23, 9).Value = 1.030831785019704
$ws.Cells.Item(23, 10).Value = 1.030850628711641
$ws.Cells.Item(23, 11).Value = 1.030785355195706
$ws.Cells.Item(23, 12).Value = 1.036685620116332
$ws.Cells.Item(23, 13).Value = 1.043714352070957
$ws.Cells.Item(23, 14).Value = 1.014307916610288
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025463634954271
$ws.Cells.Item(24, 4).Value = 1.028380564276294
$ws.Cells.Item(24, 5).Value = 1.034325248177639
$ws.Cells.Item(24, 6).Value = 1.041623159406706
$ws.Cells.Item(24, 9).Value = 1.03104346557566
$ws.Cells.Item(24, 10).Value = 1.031724897698382
$ws.Cells.Item(24, 11).Value = 1.031777432130937
$ws.Cells.Item(24, 12).Value = 1.037701164352387
$ws.Cells.Item(24, 13).Value = 1.044973714490982
$ws.Cells.Item(24, 14).Value = 1.014597381643229
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.027090410864232
$ws.Cells.Item(25, 4).Value = 1.029855180698165
$ws.Cells.Item(25, 5).Value = 1.035825956613088
$ws.Cells.Item(25, 6).Value = 1.043406657561706
$ws.Cells.Item(25, 9).Value = 1.031278773885888
$ws.Cells.Item(25, 10).Value = 1.03273616618273
$ws.Cells.Item(25, 11).Value = 1.0329265108154
$ws.Cells.Item(25, 12).Value = 1.038878404324673
$ws.Cells.Item(25, 13).Value = 1.018169663294515
$ws.Cells.Item(25, 14).Value = 1.014932116885565
